# Refresh the cryptocurrency price/volume snapshot to match the latest export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D/E value refresh (Price / Volume(1h)) ---
$ws.Range('D2').Value = '''51.915.07'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '''2.817.17'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''354.20'
$ws.Range('E5').Value = '  +6.57%  '
$ws.Range('D6').Value = '''113.51'
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('E7').Value = '  +2.42%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.599'
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('D10').Value = '''41.59'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').Value = '''0.0851'
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('D14').Value = '''7.71'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').Value = '''3.240.73'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').Value = '''2.833.23'
$ws.Range('E16').Value = '  +2.98%  '
$ws.Range('D17').Value = '''0.897'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').Value = '''51.809.84'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '''7.41'
$ws.Range('E19').Value = '  +8.32%  '
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').Value = '''13.58'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').Value = '''0.0₃0994'
$ws.Range('E22').Value = '  +2.07%  '
$ws.Range('D23').Value = '''269.97'
$ws.Range('E23').Value = '  -2.78%  '
$ws.Range('D24').Value = '''69.81'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = '''2.78'
$ws.Range('E25').Value = '  +5.13%  '
$ws.Range('D26').Value = '''26.70'
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''10.30'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('D31').Value = '''0.0458'
$ws.Range('E31').Value = '  +32.81%  '
$ws.Range('D32').Value = '''50.66'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('D33').Value = '''33.79'
$ws.Range('E33').Value = '  -3.69%  '
$ws.Range('D34').Value = '''5.82'
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').Value = '''2.08'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '''4.91'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').Value = '''3.20'
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('D40').Value = '''18.30'
$ws.Range('E40').Value = '  -4.01%  '
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''3.35'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').Value = '''2.080.18'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('E48').Value = '  +3.98%  '
$ws.Range('D49').Value = '''5.68'
$ws.Range('E49').Value = '  +2.91%  '
$ws.Range('D50').Value = '''0.938'
$ws.Range('E50').Value = '  +8.16%  '
$ws.Range('E51').Value = '  +1.30%  '

# --- Rows 41-43 got reshuffled in the source ranking (Monero dropped, EnergySwap/Stacks moved up) ---
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''23.59'
$ws.Range('E41').Value = '  +2.52%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''2.56'
$ws.Range('E42').Value = '  +4.89%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '''127.52'
$ws.Range('E43').Value = '  +0.46%  '
